# Update route for create.handlebars file.
# This updates the flashy_database_setup schema diagram:
#  - Users table: "name"/"email" -> "displayname"/"username", add setID (hasMany, blue)
#  - Sets table: swap categoryID/userID order, add flashcardID (hasMany, blue)
#  - Flashcards table: "flashNumber" -> "flash_num"
#  - Categories table: "categoryName" -> "cat_name", add setID (hasMany, blue)
#  - Add legend: GREEN = belongsTo / BLUE = hasMany

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors used in the schema (matches existing green FF00B050 "belongsTo" color already in file)
$green = 5287936    # RGB(0,176,80)  -> FF00B050 (belongsTo, already used in workbook)
$blue  = 15773696   # RGB(0,176,240) -> FF00B0F0 (hasMany, newly introduced)

# --- Legend (rows 19-20) ---
$ws.Range("A19").Value = "GREEN = belongsTo"
$ws.Range("A19").Font.Color = $green

$ws.Range("A20").Value = "BLUE = hasMany"
$ws.Range("A20").Font.Color = $blue

# --- Users table (row 2) ---
$ws.Range("B2").Value = "displayname"
$ws.Range("C2").Value = "username"
$ws.Range("E2").Value = "setID"
$ws.Range("E2").Font.Color = $blue

# --- Categories table (row 15) ---
$ws.Range("B15").Value = "cat_name"
$ws.Range("C15").Value = "setID"
$ws.Range("C15").Font.Color = $blue

# --- Flashcards table (row 11) ---
$ws.Range("B11").Value = "flash_num"

# --- Sets table (row 7): swap categoryID / userID columns ---
$ws.Range("D7").Value = "categoryID"
$ws.Range("E7").Value = "userID"
$ws.Range("F7").Value = "flashcardID"
$ws.Range("F7").Font.Color = $blue

# --- Selection moves to D5 per saved view state ---
$ws.Range("D5").Select()
